$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Earned Value")
Write-Host "Sheets:"
foreach ($s in $wb.Worksheets) {
    Write-Host $s.Name
}
